# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Updates the metrics computed for row 3 (file_name = metrics_sim_with_priors.json)
# in the metrics_sim output table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.8888888888888888
$ws.Range("H3").Value = 0.8020555767034641
$ws.Range("I3").Value = 0.03733618233618233
$ws.Range("J3").Value = 0.7777777777777778
$ws.Range("K3").Value = 111.0740740740741

$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 18
$ws.Range("S3").Value = 32
$ws.Range("T3").Value = 96
$ws.Range("U3").Value = 228
$ws.Range("V3").Value = 2598
$ws.Range("W3").Value = 2582
$ws.Range("X3").Value = 2568
$ws.Range("Y3").Value = 2504
$ws.Range("Z3").Value = 2372

$ws.Range("AF3").Value = 0.999231
$ws.Range("AG3").Value = 0.993077
$ws.Range("AH3").Value = 0.987692
$ws.Range("AI3").Value = 0.963077
$ws.Range("AJ3").Value = 0.912308
